$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.046.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.882.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.876.19'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.16%  '

$ws.Range("E8").Value = '  +0.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.547.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.896.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.264.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.12%  '

$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '467.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.738'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000160'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.046.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.863.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.15%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.140'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.61%  '

$ws.Range("E42").Value = '  +0.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.313'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000303'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '423.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '47.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.86%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.45%  '
